$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.212.79"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.247.12"
$ws.Range("E3").Value = "  -3.37%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.28"
$ws.Range("E5").Value = "  -2.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.45"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("E7").Value = "  -1.93%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -1.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.36"
$ws.Range("E10").Value = "  +2.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -3.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.94"
$ws.Range("E12").Value = "  -10.37%  "

$ws.Range("E13").Value = "  -1.78%  "

$ws.Range("D14").Value = "2.594.83"
$ws.Range("E14").Value = "  -3.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.39"
$ws.Range("E15").Value = "  +0.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.31"
$ws.Range("E16").Value = "  -3.10%  "

$ws.Range("D17").Value = "2.244.35"
$ws.Range("E17").Value = "  -3.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.727"
$ws.Range("E18").Value = "  -3.27%  "

$ws.Range("D19").Value = "40.119.77"
$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.86"
$ws.Range("E21").Value = "  -3.75%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.63"
$ws.Range("E23").Value = "  -3.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "230.96"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  -2.88%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.85"
$ws.Range("E27").Value = "  +3.32%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.07"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  +3.18%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.25"
$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.36"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.78"
$ws.Range("E32").Value = "  -0.90%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.90"
$ws.Range("E34").Value = "  -3.40%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0708"
$ws.Range("E36").Value = "  -0.60%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.65"
$ws.Range("E37").Value = "  +8.53%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0993"
$ws.Range("E39").Value = "  +0.53%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.69"
$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.942.65"
$ws.Range("E43").Value = "  +0.30%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0269"
$ws.Range("E44").Value = "  +3.30%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.10"
$ws.Range("E45").Value = "  -5.44%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.95"
$ws.Range("E46").Value = "  -2.23%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.41"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.460.96"
$ws.Range("E49").Value = "  -4.42%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.57"
$ws.Range("E50").Value = "  +2.74%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.22"
$ws.Range("E51").Value = "  -2.05%  "
